$d = $word.ActiveDocument
$paras = $d.Paragraphs

# Locate the scraped site-generator footer paragraph ("Ver no Jupiter
# Salvar em pdf Salvar em docx") that was pulled in along with the
# course page content.
$jupiterIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text -like "*Ver no Jupiter*") {
        $jupiterIndex = $i
        break
    }
}

if ($jupiterIndex -gt 0) {
    # Figure out how much to remove around it:
    #  - the blank paragraph immediately before it (if any), and
    #  - the copyright/footer paragraph ("© ... Attribution") immediately
    #    after it, if present.
    $startIndex = $jupiterIndex
    if ($jupiterIndex -gt 1) {
        $startIndex = $jupiterIndex - 1
    }

    $endIndex = $jupiterIndex
    if (($jupiterIndex + 1) -le $paras.Count -and
        ($paras.Item($jupiterIndex + 1).Range.Text -like "*Contact:*")) {
        $endIndex = $jupiterIndex + 1
    }

    $start = $paras.Item($startIndex).Range.Start
    $end = $paras.Item($endIndex).Range.End

    $r = $d.Range($start, $end)
    $r.Delete()
}
